$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price/Volume columns from Excel's automatic numeric/date
# inference so the numeric-looking text values stay as text, matching
# the original inline-string cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '31.417.41'

$ws.Range("D3").Value = '2.006.33'
$ws.Range("E3").Value = '  +7.42%  '

$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").Value = '0.7789'
$ws.Range("E5").Value = '  +64.94%  '

$ws.Range("D6").Value = '259.75'
$ws.Range("E6").Value = '  +6.77%  '

$ws.Range("D7").Value = '0.9980'
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("D8").Value = '0.3605'
$ws.Range("E8").Value = '  +25.71%  '

$ws.Range("D9").Value = '28.45'
$ws.Range("E9").Value = '  +31.45%  '

$ws.Range("D10").Value = '0.07073'
$ws.Range("E10").Value = '  +9.29%  '

$ws.Range("D11").Value = '0.8400'

$ws.Range("D12").Value = '0.08091'
$ws.Range("E12").Value = '  +4.02%  '

$ws.Range("D13").Value = '2.006.50'
$ws.Range("E13").Value = '  +7.48%  '

$ws.Range("D14").Value = '101.25'
$ws.Range("E14").Value = '  +5.38%  '

$ws.Range("D15").Value = '5.659'
$ws.Range("E15").Value = '  +10.68%  '

$ws.Range("D16").Value = '274.85'
$ws.Range("E16").Value = '  -2.52%  '

$ws.Range("D17").Value = '31.412.42'
$ws.Range("E17").Value = '  +3.69%  '

$ws.Range("D18").Value = '14.70'
$ws.Range("E18").Value = '  +13.51%  '

$ws.Range("D19").Value = '5.943'
$ws.Range("E19").Value = '  +13.51%  '

$ws.Range("D20").Value = '0.000007969'
$ws.Range("E20").Value = '  +6.98%  '

$ws.Range("D21").Value = '2.269.36'
$ws.Range("E21").Value = '  +7.66%  '

$ws.Range("D22").Value = '0.9988'
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").Value = '0.9973'
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("D24").Value = '7.215'
$ws.Range("E24").Value = '  +15.67%  '

$ws.Range("D25").Value = '10.16'
$ws.Range("E25").Value = '  +13.36%  '

$ws.Range("D26").Value = '164.13'
$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("D27").Value = '0.1479'
$ws.Range("E27").Value = '  +54.10%  '

$ws.Range("D28").Value = '20.12'
$ws.Range("E28").Value = '  +7.77%  '

$ws.Range("D29").Value = '2.396'
$ws.Range("E29").Value = '  +27.82%  '

$ws.Range("D30").Value = '1.628'
$ws.Range("E30").Value = '  +10.15%  '

$ws.Range("D31").Value = '4.631'
$ws.Range("E31").Value = '  +10.31%  '

$ws.Range("D32").Value = '1.354'
$ws.Range("E32").Value = '  +3.17%  '

$ws.Range("D33").Value = '4.394'
$ws.Range("E33").Value = '  +6.98%  '

$ws.Range("D34").Value = '0.05207'
$ws.Range("E34").Value = '  +8.86%  '

$ws.Range("D35").Value = '1.230'
$ws.Range("E35").Value = '  +10.35%  '

$ws.Range("D36").Value = '0.7602'
$ws.Range("E36").Value = '  +11.30%  '

$ws.Range("D37").Value = '2.789'
$ws.Range("E37").Value = '  +2.91%  '

$ws.Range("E38").Value = '  +7.51%  '

$ws.Range("D39").Value = '2.960'
$ws.Range("E39").Value = '  +4.19%  '

$ws.Range("D40").Value = '6.701'
$ws.Range("E40").Value = '  +8.08%  '

$ws.Range("D41").Value = '80.29'
$ws.Range("E41").Value = '  +6.53%  '

$ws.Range("D42").Value = '2.189'
$ws.Range("E42").Value = '  +14.13%  '

$ws.Range("D43").Value = '0.4741'
$ws.Range("E43").Value = '  +13.29%  '

$ws.Range("D44").Value = '0.8623'
$ws.Range("E44").Value = '  +4.91%  '

$ws.Range("D45").Value = '104.78'
$ws.Range("E45").Value = '  +4.34%  '

$ws.Range("D46").Value = '0.9991'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("D47").Value = '7.679'
$ws.Range("E47").Value = '  +9.92%  '

$ws.Range("D48").Value = '9.960'
$ws.Range("E48").Value = '  +3.94%  '

$ws.Range("D49").Value = '0.4370'
$ws.Range("E49").Value = '  +12.94%  '

$ws.Range("D50").Value = '37.12'
$ws.Range("E50").Value = '  +6.21%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '941.97'
$ws.Range("E51").Value = '  +5.63%  '

# Restore the default (unstyled) cell style now that the values are
# safely stored as text, so no stray number-format style lingers.
$ws.Range("D2:E51").Style = "Normal"
